# Insert a new data row at row 704 (pushes existing rows 704:787 down to 705:788)
# and populate it with the new weekly observation for Brócoli at Vega Modelo de Temuco.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("704:704").Insert()

$ws.Range("A704").Value = 10
$ws.Range("B704").Value = "Vega Modelo de Temuco"
$ws.Range("C704").Value = "La Araucanía"
$ws.Range("D704").Value = 45212
$ws.Range("E704").Value = 9
$ws.Range("F704").Value = 100112023
$ws.Range("G704").Value = "Brócoli"
$ws.Range("H704").Value = "Sin especificar"
$ws.Range("I704").Value = "Primera"
$ws.Range("J704").Value = 1200
$ws.Range("K704").Value = 1000
$ws.Range("L704").Value = 1000
$ws.Range("M704").Value = 1000
$ws.Range("N704").Value = "$/unidad"
$ws.Range("O704").Value = "Región del Maule"
$ws.Range("P704").Value = 1000
$ws.Range("Q704").Value = 1
$ws.Range("R704").Value = "Hortaliza"
